$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.245.72'
$ws.Range("E2").Value = '  -0.40%  '

$ws.Range("D3").Value = '3.015.23'
$ws.Range("E3").Value = '  +0.04%  '

$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '594.42'
$ws.Range("E5").Value = '  +0.98%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '148.09'
$ws.Range("E6").Value = '  +0.67%  '

$ws.Range("D8").Value = '3.016.43'
$ws.Range("E8").Value = '  +0.15%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.519'
$ws.Range("E9").Value = '  -1.67%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.38'
$ws.Range("E10").Value = '  +9.99%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.150'
$ws.Range("E11").Value = '  +0.76%  '

$ws.Range("E12").Value = '  -1.29%  '

$ws.Range("E13").Value = '  +0.89%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.44'
$ws.Range("E14").Value = '  -0.82%  '

$ws.Range("E15").Value = '  +2.36%  '

$ws.Range("D16").Value = '3.515.86'
$ws.Range("E16").Value = '  +0.08%  '

$ws.Range("B17").Value = 'Polkadot'
$ws.Range("C17").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '7.00'
$ws.Range("E17").Value = '  -1.30%  '

$ws.Range("B18").Value = 'WrappedBTC'
$ws.Range("C18").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D18").Value = '62.284.97'
$ws.Range("E18").Value = '  -0.31%  '

$ws.Range("D19").Value = '3.016.62'
$ws.Range("E19").Value = '  +0.20%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '448.00'
$ws.Range("E20").Value = '  -2.61%  '

$ws.Range("E21").Value = '  +1.38%  '

$ws.Range("E22").Value = '  -0.30%  '

$ws.Range("E23").Value = '  -0.20%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '82.32'

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '10.88'
$ws.Range("E25").Value = '  +10.28%  '

$ws.Range("E26").Value = '  +1.25%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.06'
$ws.Range("E27").Value = '  -1.53%  '

$ws.Range("E29").Value = '  +2.00%  '

$ws.Range("E30").Value = '  -0.04%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.18'
$ws.Range("E31").Value = '  +2.58%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.13'
$ws.Range("E32").Value = '  +1.42%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '27.48'
$ws.Range("E33").Value = '  -1.20%  '

$ws.Range("E34").Value = '  +0.08%  '

$ws.Range("E35").Value = '  +4.69%  '

$ws.Range("E36").Value = '  +0.20%  '

$ws.Range("E37").Value = '  +1.14%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '50.17'
$ws.Range("E38").Value = '  -0.46%  '

$ws.Range("B39").Value = 'Stacks'
$ws.Range("C39").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.06'
$ws.Range("E39").Value = '  -2.62%  '

$ws.Range("B40").Value = 'Cosmos'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '9.08'
$ws.Range("E40").Value = '  -1.52%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.97'
$ws.Range("E41").Value = '  +2.28%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.124'
$ws.Range("E42").Value = '  -0.31%  '

$ws.Range("B43").Value = 'TheGraph'
$ws.Range("C43").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.285'
$ws.Range("E43").Value = '  +6.17%  '

$ws.Range("B44").Value = 'Arweave'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '41.05'
$ws.Range("E44").Value = '  +9.52%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '394.63'
$ws.Range("E45").Value = '  -0.02%  '

$ws.Range("E46").Value = '  -1.80%  '

$ws.Range("D47").Value = '2.728.67'
$ws.Range("E47").Value = '  -0.48%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '134.78'
$ws.Range("E48").Value = '  +4.10%  '

$ws.Range("E49").Value = '  +0.10%  '

$ws.Range("E50").Value = '  -0.66%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.107'
$ws.Range("E51").Value = '  -1.65%  '
